$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Columns.Item(5).Insert()
$ws.Range("E10").Value = $false
$ws.Range("E10").Interior.ColorIndex = 0
